$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "datos actualizados" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 1 de Mayo de 2020 a las 14:22"

# Country name swaps caused by re-sorting (Republica de Macedonia / Nueva Zelanda,
# and the Sierra Leona / Etiopia / Madagascar / Guayana Francesa / Cabo Verde / Camboya block)
# plus refreshed case-count statistics for rows whose totals changed.

# Row 17
$ws.Range("B17").Value = 39791
$ws.Range("C17").Value = 475
$ws.Range("E17").Value = 34648
$ws.Range("G17").Value = 98
$ws.Range("H17").Value = 4893
# Row 21
$ws.Range("B21").Value = 25351
$ws.Range("C21").Value = 306
$ws.Range("D21").Value = 1647
$ws.Range("E21").Value = 22697
$ws.Range("F21").Value = 154
$ws.Range("G21").Value = 18
$ws.Range("H21").Value = 1007
# Row 24
$ws.Range("B24").Value = 21520
$ws.Range("C24").Value = 428
$ws.Range("E24").Value = 17862
$ws.Range("G24").Value = 67
$ws.Range("H24").Value = 2653
# Row 41
$ws.Range("D41").Value = 6729
$ws.Range("E41").Value = 2122
$ws.Range("F41").Value = 61
$ws.Range("G41").Value = 8
$ws.Range("H41").Value = 460
# Row 70
$ws.Range("B70").Value = 2085
$ws.Range("C70").Value = 9
$ws.Range("D70").Value = 1421
$ws.Range("E70").Value = 589
$ws.Range("F70").Value = 17
$ws.Range("G70").Value = 6
$ws.Range("H70").Value = 75
# Row 82
$ws.Range("A82").Value = "Republica de Macedonia"
$ws.Range("B82").Value = 1494
$ws.Range("C82").Value = 29
$ws.Range("D82").Value = 807
$ws.Range("E82").Value = 606
$ws.Range("F82").Value = 13
$ws.Range("G82").Value = 4
$ws.Range("H82").Value = 81
# Row 83
$ws.Range("A83").Value = "Nueva Zelanda"
$ws.Range("B83").Value = 1479
$ws.Range("C83").Value = 3
$ws.Range("D83").Value = 1252
$ws.Range("E83").Value = 208
$ws.Range("F83").Value = 1
$ws.Range("H83").Value = 19
# Row 103
$ws.Range("B103").Value = 671
$ws.Range("C103").Value = 8
$ws.Range("E103").Value = 507
# Row 142
$ws.Range("A142").Value = "Sierra Leona"
$ws.Range("B142").Value = 136
$ws.Range("C142").Value = 12
$ws.Range("D142").Value = 21
$ws.Range("E142").Value = 108
$ws.Range("H142").Value = 7
# Row 143
$ws.Range("A143").Value = "Etiopia"
$ws.Range("B143").Value = 133
$ws.Range("C143").Value = 2
$ws.Range("D143").Value = 66
$ws.Range("E143").Value = 64
$ws.Range("F143").Value = 0
$ws.Range("H143").Value = 3
# Row 144
$ws.Range("A144").Value = "Madagascar"
$ws.Range("B144").Value = 132
$ws.Range("C144").Value = 4
$ws.Range("D144").Value = 94
$ws.Range("E144").Value = 38
$ws.Range("F144").Value = 1
$ws.Range("H144").Value = 0
# Row 145
$ws.Range("A145").Value = "Guayana Francesa"
$ws.Range("B145").Value = 128
$ws.Range("C145").Value = 2
$ws.Range("D145").Value = 98
$ws.Range("E145").Value = 29
$ws.Range("F145").Value = 2
$ws.Range("H145").Value = 1
# Row 146
$ws.Range("A146").Value = "Cabo Verde"
$ws.Range("C146").Value = 1
$ws.Range("D146").Value = 4
$ws.Range("E146").Value = 117
$ws.Range("F146").Value = 0
$ws.Range("H146").Value = 1
# Row 147
$ws.Range("A147").Value = "Camboya"
$ws.Range("B147").Value = 122
$ws.Range("D147").Value = 119
$ws.Range("E147").Value = 3
$ws.Range("F147").Value = 1
$ws.Range("H147").Value = 0
